# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 867
$ws1.Range("F7").Value  = 453
$ws1.Range("F9").Value  = 2121
$ws1.Range("F10").Value = 607
$ws1.Range("F11").Value = 272
$ws1.Range("F13").Value = 1022
$ws1.Range("F14").Value = 166
$ws1.Range("F15").Value = 2166
$ws1.Range("F16").Value = 617
$ws1.Range("F17").Value = 11399
$ws1.Range("F18").Value = 1160
$ws1.Range("F20").Value = 116
$ws1.Range("F21").Value = 132
$ws1.Range("F26").Value = 2

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 14
$ws2.Range("F10").Value = 13
$ws2.Range("F11").Value = 76
$ws2.Range("F12").Value = 50
$ws2.Range("F17").Value = 33
$ws2.Range("F19").Value = 2

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 5670
$ws3.Range("F4").Value = 448

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 5670
$ws4.Range("F5").Value  = 448
$ws4.Range("F10").Value = 867
$ws4.Range("F12").Value = 453
$ws4.Range("F13").Value = 14
$ws4.Range("F14").Value = 2121
$ws4.Range("F15").Value = 607
$ws4.Range("F16").Value = 272
$ws4.Range("F20").Value = 1022
$ws4.Range("F22").Value = 166
$ws4.Range("F24").Value = 13
$ws4.Range("F25").Value = 2166
$ws4.Range("F26").Value = 617
$ws4.Range("F27").Value = 11399
$ws4.Range("F28").Value = 76
$ws4.Range("F29").Value = 50
$ws4.Range("F30").Value = 1160
$ws4.Range("F32").Value = 116
$ws4.Range("F33").Value = 132
$ws4.Range("F40").Value = 33
$ws4.Range("F42").Value = 2
$ws4.Range("F49").Value = 2
